# REPORT.xlsx - "To Do" sheet updates
# store cache with static key , description in getList & acceptor and submitor in admin web

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do")

# Rename existing task text: "post submit challenge" -> "post submit challenge Test"
$ws.Range("A27").Value = "post submit challenge Test"

# Mark the "acceptor"/"submitor" admin-web rows (24 & 25) with a "++" flag in column B.
# Leading apostrophe forces text (quote-prefix) entry, matching the workbook's
# existing "++" / text-flag style used elsewhere in column B.
$ws.Range("B24").Value = "'++"
$ws.Range("B25").Value = "'++"

# The renamed "post submit challenge Test" task (row 27) is now Done.
$ws.Range("B27").Value = "Done"

# Two new tasks appended to the list.
$ws.Range("A28").Value = "add discription when sharing"
$ws.Range("A29").Value = "show description in getList"

# Leave selection on the last entered cell, matching the saved view state.
$ws.Range("A29").Select()
